# Added manual REX hours
# Update the hour values on the "prepa ORTEMS" sheet (row 2) to reflect
# manually entered REX hours (scaled by a factor of 1.57815 vs. the
# previous computed values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prepa ORTEMS")

$ws.Range("F2").Value = 87.15000000000001
$ws.Range("G2").Value = 14.9625
$ws.Range("H2").Value = 134.5575
$ws.Range("J2").Value = 6.300000000000001
$ws.Range("L2").Value = 80.85000000000001
$ws.Range("M2").Value = 175.9695
$ws.Range("N2").Value = 93.87
$ws.Range("O2").Value = 665.7105
$ws.Range("P2").Value = 4.2
$ws.Range("U2").Value = 105.105
$ws.Range("V2").Value = 12.915
$ws.Range("W2").Value = 15.015
$ws.Range("X2").Value = 25.83
$ws.Range("Y2").Value = 42.315
$ws.Range("Z2").Value = 29.4
$ws.Range("AD2").Value = 3.675
$ws.Range("AF2").Value = 5.775
$ws.Range("AM2").Value = 27.3
$ws.Range("AS2").Value = 24.15
$ws.Range("AT2").Value = 2.1
$ws.Range("AU2").Value = 4.2
$ws.Range("AV2").Value = 2.1
$ws.Range("AW2").Value = 6.300000000000001
$ws.Range("AX2").Value = 8.4
